$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.764.22'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.644.72'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.47%  '
$ws.Range('D5').Value = '216.72'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').Value = '19.13'
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '1.649.56'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').Value = '64.64'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').Value = '26.753.66'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = '0.0₃0735'
$ws.Range('E17').Value = '  -1.86%  '
$ws.Range('D18').Value = '213.96'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D20').Value = '4.39'
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('E21').Value = '  +12.55%  '
$ws.Range('D22').Value = '6.24'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').Value = '9.34'
$ws.Range('E23').Value = '  -2.18%  '
$ws.Range('D24').Value = '146.04'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').Value = '15.63'
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').Value = '3.37'
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('D32').Value = '3.01'
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('D33').Value = '1.290.73'
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').Value = '2.44'
$ws.Range('E35').Value = '  +1.43%  '
$ws.Range('E36').Value = '  -3.07%  '
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('D38').Value = '0.818'
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '0.804'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.23'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('D43').Value = '1.789.31'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').Value = '61.52'
$ws.Range('E44').Value = '  +3.55%  '
$ws.Range('D45').Value = '91.69'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').Value = '1.61'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').Value = '0.0₆0104'
$ws.Range('E47').Value = '  -1.63%  '
$ws.Range('D48').Value = '0.0523'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('D49').Value = '7.66'
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('E50').Value = '  -0.04%  '
$ws.Range('E51').Value = '  +0.07%  '
